# Updated capital structure database
# Recomputed financial ratios for Saint Lucia / Investments & Asset Management rows (2-5)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Row 2
$ws.Range("G2").Value = 0.03248587570621469
$ws.Range("H2").Value = 0.03248587570621469
$ws.Range("I2").Value = 1.297316384180791
$ws.Range("J2").Value = 1.29353412067006
$ws.Range("K2").Value = -12.77
$ws.Range("L2").Value = 1.503060263653484
$ws.Range("M2").Value = 0.731
$ws.Range("N2").Value = 0.00753841394245643
$ws.Range("O2").Value = -0.0572435395458105
$ws.Range("P2").Value = 0.731
$ws.Range("Q2").Value = 0.00753841394245643
$ws.Range("R2").Value = -0.0572435395458105
$ws.Range("U2").Value = 4.34
$ws.Range("V2").Value = 0.04475611013715582
$ws.Range("W2").Value = -0.04960317460317461
$ws.Range("X2").Value = 0.06175660787539299
$ws.Range("Y2").Value = -0.1113597824785676
$ws.Range("Z2").Value = -0.0470922111611201
$ws.Range("AA2").Value = -0.04583476566976052
$ws.Range("AB2").Value = 0.05893305391276064
$ws.Range("AC2").Value = -0.1035452081556093
$ws.Range("AD2").Value = 22.22
$ws.Range("AF2").Value = 22.22
$ws.Range("AG2").Value = 17.88
$ws.Range("AH2").Value = 0.1864250356573538
$ws.Range("AI2").Value = 0.1773768659695058
$ws.Range("AJ2").Value = 0.1556813234653896
$ws.Range("AK2").Value = 0.1478541304887125
$ws.Range("AL2").Value = 1.861
$ws.Range("AM2").Value = 1.861
$ws.Range("AO2").Value = -5.922622246104244
$ws.Range("AQ2").Value = -5.922622246104244
# Row 3
$ws.Range("I3").Value = 0.6872770511296076
$ws.Range("J3").Value = 0.6812658815861999
$ws.Range("K3").Value = 0.68
$ws.Range("L3").Value = 0.8085612366230679
$ws.Range("M3").Value = 0.205
$ws.Range("N3").Value = 0.02638352638352638
$ws.Range("O3").Value = 0.3014705882352941
$ws.Range("P3").Value = 0.205
$ws.Range("Q3").Value = 0.02638352638352638
$ws.Range("R3").Value = 0.3014705882352941
$ws.Range("U3").Value = 0.008
$ws.Range("V3").Value = 0.00102960102960103
$ws.Range("W3").Value = 0.07497243660418963
$ws.Range("X3").Value = 0.06857950216002676
$ws.Range("Y3").Value = 0.006392934444162876
$ws.Range("Z3").Value = 0.07344978165938865
$ws.Range("AA3").Value = 0.0500388302544973
$ws.Range("AB3").Value = 0.06212199305859818
$ws.Range("AC3").Value = -0.01208316280410087
$ws.Range("AD3").Value = 3.21
$ws.Range("AF3").Value = 3.21
$ws.Range("AG3").Value = 3.202
$ws.Range("AH3").Value = 0.2923497267759563
$ws.Range("AI3").Value = 0.2706576728499157
$ws.Range("AJ3").Value = 0.2918337586584032
$ws.Range("AK3").Value = 0.2701653729328383
$ws.Range("AL3").Value = 0.185
$ws.Range("AM3").Value = 0.185
$ws.Range("AO3").Value = 3.124324324324324
$ws.Range("AQ3").Value = 3.124324324324324
# Row 4
$ws.Range("G4").Value = -0
$ws.Range("H4").Value = -0
$ws.Range("I4").Value = 2.044293015332198
$ws.Range("J4").Value = 2.044293015332198
$ws.Range("K4").Value = -1.25
$ws.Range("L4").Value = 2.129471890971039
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 0.002
$ws.Range("V4").Value = 0.00009569377990430622
$ws.Range("W4").Value = -0.04960317460317461
$ws.Range("X4").Value = 0.05959463105643872
$ws.Range("Y4").Value = -0.1091978056596133
$ws.Range("Z4").Value = -0.02242083954012452
$ws.Range("AA4").Value = -0.04583476566976052
$ws.Range("AB4").Value = 0.05771044248584881
$ws.Range("AC4").Value = -0.1035452081556093
$ws.Range("AD4").Value = 3.51
$ws.Range("AF4").Value = 3.51
$ws.Range("AG4").Value = 3.508
$ws.Range("AH4").Value = 0.1437935272429332
$ws.Range("AI4").Value = 0.1275899672846238
$ws.Range("AJ4").Value = 0.1437233693870862
$ws.Range("AK4").Value = 0.1275265377344773
$ws.Range("AL4").Value = 0.326
$ws.Range("AM4").Value = 0.326
$ws.Range("AO4").Value = -3.680981595092024
$ws.Range("AQ4").Value = -3.680981595092024
# Row 5
$ws.Range("G5").Value = 0.03154285714285714
$ws.Range("H5").Value = 0.03154285714285714
$ws.Range("I5").Value = 1.188571428571429
$ws.Range("J5").Value = 1.188571428571429
$ws.Range("K5").Value = -12.2
$ws.Range("L5").Value = 1.394285714285714
$ws.Range("M5").Value = 0.526
$ws.Range("N5").Value = 0.007701317715959005
$ws.Range("O5").Value = -0.04311475409836066
$ws.Range("P5").Value = 0.526
$ws.Range("Q5").Value = 0.007701317715959005
$ws.Range("R5").Value = -0.04311475409836066
$ws.Range("U5").Value = 4.33
$ws.Range("V5").Value = 0.06339677891654466
$ws.Range("W5").Value = -0.09501557632398754
$ws.Range("X5").Value = 0.06175660787539299
$ws.Range("Y5").Value = -0.1567721841993805
$ws.Range("Z5").Value = -0.06128266365973062
$ws.Range("AA5").Value = -0.07283882309270839
$ws.Range("AB5").Value = 0.05893305391276064
$ws.Range("AC5").Value = -0.131771877005469
$ws.Range("AD5").Value = 15.5
$ws.Range("AF5").Value = 15.5
$ws.Range("AG5").Value = 11.17
$ws.Range("AH5").Value = 0.184964200477327
$ws.Range("AI5").Value = 0.1804423748544819
$ws.Range("AJ5").Value = 0.1405561847237951
$ws.Range("AK5").Value = 0.1369375996076989
$ws.Range("AL5").Value = 1.35
$ws.Range("AM5").Value = 1.35
$ws.Range("AO5").Value = -7.703703703703703
$ws.Range("AQ5").Value = -7.703703703703703
